# Update Work Week and Social Spending
#
# Rwanda "GDP per Capita" time series (sheet "Data"):
#   - refresh the per-year "Data" values for the existing years
#     1950-2008 (rows 2-60) with the newer revised figures, and
#   - extend the series with 8 new years, 2009-2016 (rows 61-68).
#
# The "Data" column stores the figures as text (as in the source file),
# so each value is written with a leading apostrophe to force text type
# instead of a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$years = @(1950, 1951, 1952, 1953, 1954, 1955, 1956, 1957, 1958, 1959, 1960, 1961, 1962, 1963, 1964, 1965, 1966, 1967, 1968, 1969, 1970, 1971, 1972, 1973, 1974, 1975, 1976, 1977, 1978, 1979, 1980, 1981, 1982, 1983, 1984, 1985, 1986, 1987, 1988, 1989, 1990, 1991, 1992, 1993, 1994, 1995, 1996, 1997, 1998, 1999, 2000, 2001, 2002, 2003, 2004, 2005, 2006, 2007, 2008, 2009, 2010, 2011, 2012, 2013, 2014, 2015, 2016)

$values = @("'1020", "'1058", "'1071", "'1090", "'1127", "'1138", "'1156", "'1173", "'1178", "'1216", "'1224", "'1167", "'1296", "'1140", "'980", "'1023", "'1065", "'1109", "'1154", "'1243", "'1337", "'1315", "'1282", "'1288", "'1262", "'1199", "'1387", "'1409", "'1495", "'1584", "'1672", "'1706", "'1680", "'1722", "'1597", "'1621", "'1642", "'1586", "'1538", "'1476", "'1398", "'1282.25940200987", "'1301.68900792678", "'1110.92085707389", "'707.297099996548", "'1049.23853407576", "'976.407910612693", "'932.005892157274", "'950.389473688911", "'975.835922473427", "'1023.61822713066", "'1067.2347821601", "'1154.10609348501", "'1130.04140421816", "'1159.02810839957", "'1210.51054312556", "'1261.76091707443", "'1291.19332531318", "'1364.74501289332", "'1384.2588512137", "'1416.39230711648", "'1449", "'1531", "'1554", "'1615", "'1693", "'1758")

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = 2 + $i
    $year = $years[$i]
    $val = $values[$i]

    # Rows beyond the previous last row (60, year 2008) are brand new —
    # fill in the Country Code / Country Name / Indicator / Year columns too.
    if ($row -gt 60) {
        $ws.Cells.Item($row, 1).Value = 646
        $ws.Cells.Item($row, 2).Value = "Rwanda"
        $ws.Cells.Item($row, 3).Value = "GDP per Capita"
        $ws.Cells.Item($row, 4).Value = $year
    }

    $ws.Cells.Item($row, 5).Value = $val
}
